$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: new client record (Luis Paniagua Sabines) -------------------
# Previously D20 held an empty "clien" placeholder string (shared string
# now renamed to "paniagua") and F20 held that placeholder text. Replace
# the whole row with the real record.
$ws.Range("A20").Value = "Luis"
$ws.Range("B20").Value = "paniagua"
$ws.Range("C20").Value = "Sabines"
$ws.Range("D20").Value = "lpaniagua@hotmail.com"
$ws.Range("E20").Value = 5534236787
$ws.Range("F20").Value = 5543231278

# Give D20 the same "Hipervinculo" (Hyperlink) cell style the other mail
# cells in column D use, then attach the actual hyperlink.
$ws.Range("D20").Style = "Hipervínculo"
$ws.Hyperlinks.Add($ws.Range("D20"), "mailto:lpaniagua@hotmail.com", "", "", "lpaniagua@hotmail.com")
$ws.Range("D20").Style = "Hipervínculo"

# --- Selection moves to E13 ------------------------------------------------
$ws.Range("E13").Select() | Out-Null
